$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D look numeric (e.g. "59.282.37", "123.00") but must
# stay literal text, exactly like the source inline strings. A leading
# apostrophe forces Excel to store them as text instead of coercing to a number.

$ws.Range("D2").Value = "'59.282.37"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "'2.514.79"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'534.79"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").Value = "'139.35"
$ws.Range("E6").Value = "  -3.86%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'2.518.25"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  +1.50%  "
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "'0.356"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'2.962.04"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'23.47"
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "'59.196.19"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'0.0000141"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "'2.516.86"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'11.14"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").Value = "'4.31"
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'325.07"
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "'5.80"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("D24").Value = "'63.80"
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("D25").Value = "'0.429"
$ws.Range("E25").Value = "  -1.48%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").Value = "'7.84"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "'6.89"
$ws.Range("E29").Value = "  +1.78%  "
$ws.Range("D30").Value = "'0.0₃0776"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").Value = "'1.79"
$ws.Range("E31").Value = "  -2.46%  "
$ws.Range("D32").Value = "'165.11"
$ws.Range("E32").Value = "  +5.19%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").Value = "'1.46"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("E35").Value = "  -8.61%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").Value = "'4.26"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").Value = "'1.58"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").Value = "'36.89"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("D40").Value = "'3.70"
$ws.Range("E40").Value = "  -0.58%  "
$ws.Range("D41").Value = "'0.817"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").Value = "'5.25"
$ws.Range("E42").Value = "  -6.92%  "
$ws.Range("D43").Value = "'279.27"
$ws.Range("E43").Value = "  -5.96%  "
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").Value = "'10.87"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").Value = "'0.0933"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").Value = "'123.00"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "'17.80"
$ws.Range("E51").Value = "  -2.64%  "
